$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the SERVICE_NAME column entirely (service assignment removed),
# shifting everything from column E onward left by one.
$ws.Columns("D").Delete()

# Fix the DATA_OF_BIRTH typo -> DATE_OF_BIRTH (now in column F after the shift)
$ws.Range("F1").Value = "DATE_OF_BIRTH"

# Shrink row 2 to its natural (non-custom) height
$ws.Rows(2).RowHeight = 12.75

# Extend the wrap-text formatted column (now E) down through rows 4-5
$ws.Range("E2").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E5").PasteSpecial(-4122)
$ws.Rows(4).RowHeight = 15.75
$ws.Rows(5).RowHeight = 15.75

# Update the active selection to match the saved view
$ws.Range("E9").Select()
